# Update countries & provincias Spain
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the "last updated" timestamp string (row 1) ---
$ws.Range("A1").Value = "Datos actualizados a 30 de Mayo de 2020 a las 23:05"

# --- Row 4: Estados Unidos - updated counters ---
$ws.Range("B4").Value = 1811767
$ws.Range("C4").Value = 18237
$ws.Range("D4").Value = 528172
$ws.Range("E4").Value = 1178240
$ws.Range("G4").Value = 813
$ws.Range("H4").Value = 105355

# --- Row 11: Alemania - updated counters ---
$ws.Range("B11").Value = 183281
$ws.Range("C11").Value = 262
$ws.Range("E11").Value = 9781
$ws.Range("G11").Value = 6
$ws.Range("H11").Value = 8600

# --- Rows 14/15: Iran and Peru swap order (Peru overtakes Iran) ---
# Row 14 becomes Peru with new, higher counters
$ws.Range("A14").Value = "Peru"
$ws.Range("B14").Value = 155671
$ws.Range("C14").Value = 7386
$ws.Range("D14").Value = 66447
$ws.Range("E14").Value = 84853
$ws.Range("G14").Value = 141
$ws.Range("H14").Value = 4371

# Row 15 becomes Iran, carrying the counters Iran previously had on row 14
$ws.Range("A15").Value = "Iran"
$ws.Range("B15").Value = 148950
$ws.Range("C15").Value = 2282
$ws.Range("D15").Value = 116827
$ws.Range("E15").Value = 24389
$ws.Range("G15").Value = 57
$ws.Range("H15").Value = 7734

# --- Row 53: Barein - updated counters ---
$ws.Range("E53").Value = 4950
$ws.Range("G53").Value = 2
$ws.Range("H53").Value = 17

# --- Rows 98/99: Lituania and Maldivas swap order (Maldivas overtakes Lituania) ---
# Row 98 becomes Maldivas with new, higher counters
$ws.Range("A98").Value = "Maldivas"
$ws.Range("B98").Value = 1672
$ws.Range("C98").Value = 81
$ws.Range("D98").Value = 393
$ws.Range("E98").Value = 1274
$ws.Range("G98").Value = 0
$ws.Range("H98").Value = 5

# Row 99 becomes Lituania, carrying the counters Lituania previously had on row 98
$ws.Range("A99").Value = "Lituania"
$ws.Range("B99").Value = 1670
$ws.Range("C99").Value = 8
$ws.Range("D99").Value = 1229
$ws.Range("E99").Value = 371
$ws.Range("G99").Value = 2
$ws.Range("H99").Value = 70

# --- Row 121: Niger - updated counters ---
$ws.Range("B121").Value = 956
$ws.Range("C121").Value = 1
$ws.Range("D121").Value = 818
$ws.Range("E121").Value = 74
